# Update "想去人数" (F column) figures on the exhibition/event sheets to
# reflect the newly generated output (gh-pages rebuild at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 229
$ws1.Range("F3").Value  = 1050
$ws1.Range("F5").Value  = 13810
$ws1.Range("F7").Value  = 248
$ws1.Range("F8").Value  = 1776
$ws1.Range("F9").Value  = 169
$ws1.Range("F10").Value = 133
$ws1.Range("F12").Value = 49
$ws1.Range("F13").Value = 529
$ws1.Range("F15").Value = 1
$ws1.Range("F16").Value = 13852
$ws1.Range("F17").Value = 365
$ws1.Range("F19").Value = 14925
$ws1.Range("F21").Value = 8222
$ws1.Range("F22").Value = 273
$ws1.Range("F27").Value = 162
$ws1.Range("F31").Value = 3
$ws1.Range("F38").Value = 210
$ws1.Range("F39").Value = 389
$ws1.Range("F41").Value = 5067

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 48

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 229
$ws4.Range("F3").Value  = 1050
$ws4.Range("F5").Value  = 13810
$ws4.Range("F7").Value  = 248
$ws4.Range("F8").Value  = 1776
$ws4.Range("F9").Value  = 169
$ws4.Range("F10").Value = 133
$ws4.Range("F12").Value = 49
$ws4.Range("F13").Value = 529
$ws4.Range("F15").Value = 1
$ws4.Range("F16").Value = 13852
$ws4.Range("F17").Value = 365
$ws4.Range("F19").Value = 14925
$ws4.Range("F21").Value = 8222
$ws4.Range("F22").Value = 273
$ws4.Range("F27").Value = 162
$ws4.Range("F31").Value = 3
$ws4.Range("F34").Value = 48
$ws4.Range("F40").Value = 210
$ws4.Range("F41").Value = 389
$ws4.Range("F43").Value = 5067
